$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").ClearContents()
$ws.Range("E2").ClearContents()
